$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "hội chiếu tại Thiên Di" (converging-stars-at-Travel-palace) pattern
# strings appended to the lookup table in columns A and B.
$values = @(
    'Sát Phá Tham hội chiếu tại Thiên Di',
    'Tử Phủ Vũ Tướng hội chiếu tại Thiên Di',
    'Cơ Nguyệt Đồng Lương hội chiếu tại Thiên Di',
    'Cự Nhật hội chiếu tại Thiên Di',
    'Kình Đà hội chiếu tại Thiên Di',
    'Xương Khúc hội chiếu tại Thiên Di',
    'Hoả Linh hội chiếu tại Thiên Di',
    'Không Kiếp hội chiếu tại Thiên Di',
    'Quang Quý hội chiếu tại Thiên Di',
    'Tả Hữu hội chiếu tại Thiên Di',
    'Song Hao hội chiếu tại Thiên Di',
    'Tang Hổ hội chiếu tại Thiên Di',
    'Khốc Hư hội chiếu tại Thiên Di',
    'Hình Riêu hội chiếu tại Thiên Di',
    'Thai Toạ hội chiếu tại Thiên Di',
    'Đào Hồng hội chiếu tại Thiên Di',
    'Ấn Phù hội chiếu tại Thiên Di',
    'Song Hao Quyền Lộc Kiếp Hoả hội chiếu tại Thiên Di',
    'Tử Phủ Vũ Tướng Xương Khúc Khôi Việt Tả Hữu Khoa Quyền Lộc Long hội chiếu tại Thiên Di',
    'Tử Phủ Vũ Tướng Tả Hữu Khoa Quyền Lộc Long Phượng hội chiếu tại Thiên Di',
    'Tử Khúc Phá Dương Đà hội chiếu tại Thiên Di',
    'Cơ Nguyệt Đồng Lương gặp Xương Khúc Tả Hữu hội chiếu tại Thiên Di',
    'Cơ Nguyệt Đồng Lương Khoa Tả Hữu Quang Quý Quan Phúc hội chiếu tại Thiên Di',
    'Sát Quyền hội chiếu tại Thiên Di',
    'Lộc Mã hội chiếu tại Thiên Di',
    'Kiếp Hư Hao Quyền hội chiếu tại Thiên Di',
    'Tuế Hổ Phù Xương Khúc hội chiếu tại Thiên Di',
    'Xương Khúc Tấu Long Phượng hội chiếu tại Thiên Di',
    'Đào Hồng Riêu Tấu Cơ Vũ hội chiếu tại Thiên Di',
    'Binh Hình Tướng Ấn hội chiếu tại Thiên Di',
    'Hổ Tấu hội chiếu tại Thiên Di',
    'Hình Riêu Y hội chiếu tại Thiên Di',
    'Mã Hỏa Linh hội chiếu tại Thiên Di',
    'Thai Tọa Hồng Đào hội chiếu tại Thiên Di',
    'Tả Hữu Không Kiếp hội chiếu tại Thiên Di',
    'Tả Hữu Binh Tướng hội chiếu tại Thiên Di',
    'Đào Quyền hội chiếu tại Thiên Di',
    'Đào Hồng Tả Cái Hữu hội chiếu tại Thiên Di',
    'Quan Phúc Quang Tấu hội chiếu tại Thiên Di',
    'Đào Hồng Xương Khúc hội chiếu tại Thiên Di',
    'Đào Tử Phủ hội chiếu tại Thiên Di',
    'Hổ Kình Sát hội chiếu tại Thiên Di',
    'Hổ Tang Kiếp hội chiếu tại Thiên Di',
    'Hổ Tang Không Kiếp hội chiếu tại Thiên Di',
    'Âm Dương Lương hội chiếu tại Thiên Di',
    'Cơ Lương Gia Hội hội chiếu tại Thiên Di',
    'Nhật Chiếu Lôi Môn hội chiếu tại Thiên Di',
    'Tả Hữu Xương Khúc hội chiếu tại Thiên Di',
    'Tham Linh Triều Viên hội chiếu tại Thiên Di',
    'Cự Hỏa Linh hội chiếu tại Thiên Di',
    'Hồng Đào Kỵ hội chiếu tại Thiên Di',
    'Cự Đồng Hình hội chiếu tại Thiên Di',
    'Lương Phá hội chiếu tại Thiên Di',
    'Khôi Việt hội chiếu tại Thiên Di',
    'Tham Vũ Hỏa hội chiếu tại Thiên Di',
    'Cự Kỵ hội chiếu tại Thiên Di',
    'Kình Đà Hỏa hội chiếu tại Thiên Di',
    'Kình Đà Không Kiếp hội chiếu tại Thiên Di',
    'Đào Hồng Hỷ hội chiếu tại Thiên Di',
    'Đào Hồng Kiếp Sát hội chiếu tại Thiên Di',
    'Đào Hồng Kỵ hội chiếu tại Thiên Di',
    'Đào Hồng Tả Phù Hữu Bật hội chiếu tại Thiên Di',
    'Đào Hồng Tả Phù Hữu Bật Khoa Quyền Lộc hội chiếu tại Thiên Di',
    'Đào Hồng Tả Phù Hữu Bật Khoa Quyền Lộc Long Trì Phượng Các hội chiếu tại Thiên Di',
    'Lưu Hà Kiếp Sát hội chiếu tại Thiên Di',
    'Phục Không Kiếp hội chiếu tại Thiên Di',
    'Xương Khúc Khôi Việt hội chiếu tại Thiên Di',
    'Khoa Quyền hội chiếu tại Thiên Di',
    'Lộc Quyền hội chiếu tại Thiên Di',
    'Tướng Binh Đào Hồng hội chiếu tại Thiên Di',
    'Lương Khốc Tuế hội chiếu tại Thiên Di',
    'Tướng Binh hội chiếu tại Thiên Di',
    'Cự Tang hội chiếu tại Thiên Di',
    'Cự Hỏa hội chiếu tại Thiên Di',
    'Phá Hình Kỵ hội chiếu tại Thiên Di',
    'Tang Trực Tuế hội chiếu tại Thiên Di',
    'Tuế Xương Khúc hội chiếu tại Thiên Di',
    'Tử Vi Tả Hữu hội chiếu tại Thiên Di',
    'Tử Vi Tang Tả Hữu hội chiếu tại Thiên Di'
)

$startRow = 4214
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$endRow = $startRow + $values.Length - 1

# Match the saved view state: selection on the newly-added block and the
# window scrolled so that block is visible.
[void]$ws.Range("B$startRow`:B$endRow").Select()
$excel.ActiveWindow.ScrollRow = 4185
$excel.ActiveWindow.ScrollColumn = 1
